# "Setting up run modes for Test Suites"
#
# - Rename the 3rd sheet ("Лист3") to "test_suite"
# - Populate it with a TCID / Runmode table
# - Give column A a wider, custom width
# - Make the new sheet the active tab, with B5 selected (as left by the author)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

$ws.Name = "test_suite"

# Fill the cells in the same order the original author typed them in, so the
# shared-string table ends up in the same sequence (Y before BankManagerLoginTest).
$ws.Range("A1").Value = "TCID"
$ws.Range("B1").Value = "Runmode"

$ws.Range("B2").Value = "Y"
$ws.Range("A2").Value = "BankManagerLoginTest"

$ws.Range("A3").Value = "AddCustomerTest"
$ws.Range("B3").Value = "Y"

$ws.Range("A4").Value = "OpenAccountTest"
$ws.Range("B4").Value = "N"

# Widen column A to fit the longest test-case name.
$ws.Columns.Item(1).ColumnWidth = 23.28

# Activate the sheet (moves tabSelected from sheet1 to this sheet and sets
# the workbook's activeTab) and leave the selection where the author left it.
$ws.Activate()
$ws.Range("B5").Select() | Out-Null
